$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank rows before row 211; existing rows 211-272 shift down to 213-274.
$ws.Rows.Item(211).Insert()
$ws.Rows.Item(211).Insert()

# New row 211 data
$ws.Cells.Item(211, 1).Value = 3
$ws.Cells.Item(211, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(211, 3).Value = "Coquimbo"
$ws.Cells.Item(211, 4).Value = 44900
$ws.Cells.Item(211, 5).Value = 5
$ws.Cells.Item(211, 6).Value = "Fruta"
$ws.Cells.Item(211, 7).Value = 100101
$ws.Cells.Item(211, 8).Value = "Berries"
$ws.Cells.Item(211, 9).Value = 100101001
$ws.Cells.Item(211, 10).Value = "Arándano (blue)"
$ws.Cells.Item(211, 11).Value = "Sin especificar"
$ws.Cells.Item(211, 12).Value = "Primera"
$ws.Cells.Item(211, 13).Value = 160
$ws.Cells.Item(211, 14).Value = 5800
$ws.Cells.Item(211, 15).Value = 6000
$ws.Cells.Item(211, 16).Value = 5894
$ws.Cells.Item(211, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(211, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(211, 19).Value = 2947
$ws.Cells.Item(211, 20).Value = 2

# New row 212 data
$ws.Cells.Item(212, 1).Value = 3
$ws.Cells.Item(212, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(212, 3).Value = "Coquimbo"
$ws.Cells.Item(212, 4).Value = 44900
$ws.Cells.Item(212, 5).Value = 5
$ws.Cells.Item(212, 6).Value = "Fruta"
$ws.Cells.Item(212, 7).Value = 100101
$ws.Cells.Item(212, 8).Value = "Berries"
$ws.Cells.Item(212, 9).Value = 100101001
$ws.Cells.Item(212, 10).Value = "Arándano (blue)"
$ws.Cells.Item(212, 11).Value = "Sin especificar"
$ws.Cells.Item(212, 12).Value = "Segunda"
$ws.Cells.Item(212, 13).Value = 98
$ws.Cells.Item(212, 14).Value = 4800
$ws.Cells.Item(212, 15).Value = 5000
$ws.Cells.Item(212, 16).Value = 4902
$ws.Cells.Item(212, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(212, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(212, 19).Value = 2451
$ws.Cells.Item(212, 20).Value = 2
